# Append the new resale-data row (row 69) reported 2025-02-14 08:59:10
# to the "CityResaleNum" sheet, matching the commit
# "Realestate Update resale numbers 2025-02-14 08:59".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

# Force the Date/Time/Week columns to be stored as plain text so Excel
# doesn't auto-convert values like "2025-02-14" or "06" into a date
# serial number / numeric 6 (same text formatting already used by the
# existing rows in columns A, B and D).
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("B" + $row).NumberFormat = "@"
$ws.Range("D" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-02-14"
$ws.Range("B" + $row).Value = "08:59:10"
$ws.Range("C" + $row).Value = "Friday"
$ws.Range("D" + $row).Value = "06"
$ws.Range("E" + $row).Value = 120951
$ws.Range("F" + $row).Value = 142386
$ws.Range("G" + $row).Value = 169950
$ws.Range("H" + $row).Value = 159018
$ws.Range("I" + $row).Value = -1
$ws.Range("J" + $row).Value = 144725
$ws.Range("K" + $row).Value = -1
$ws.Range("L" + $row).Value = -1
$ws.Range("M" + $row).Value = 192094
$ws.Range("N" + $row).Value = 115197
$ws.Range("O" + $row).Value = 44909
$ws.Range("P" + $row).Value = 28629
$ws.Range("Q" + $row).Value = 65347
$ws.Range("R" + $row).Value = -1
$ws.Range("S" + $row).Value = 44583
$ws.Range("T" + $row).Value = -1
